# Generate Report for Handoff
# The "2f7780e8-e76e-469e-91b9-bc99e9da7c12" file has moved from
# "Handed back: in sync with en-US" to "Ready for handoff", with updated
# handoff timestamps and a new error-detail message reporting that the
# handback file is out of date.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the 2f7780e8... file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-05 20:57:08"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ff8123df3cee00cc76399934ff8668ffb0c33842/e2e/2f7780e8-e76e-469e-91b9-bc99e9da7c12.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ba82167f7106940fa0efd4c27f68f8b62aaf919a/e2e/2f7780e8-e76e-469e-91b9-bc99e9da7c12.md."

# Column P auto-widened to fit the new Error Detail text (stored OOXML width
# of 40 == ColumnWidth 40 minus the engine's 5/6-char padding offset).
$errorColWidth = 40 - (5 / 6)

# --- zh-cn sheet: row 3 is the 2f7780e8... file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-09-05 20:56:58"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = $errorColWidth

# --- de-de sheet: row 3 is the 2f7780e8... file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-09-05 20:57:08"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = $errorColWidth
